# Remove the obsolete year rows (2000, 2005-2009) from the table, leaving
# only 2010-2013 in place. Deleting rows 2 through 7 shifts the remaining
# rows (old 8-11, for 2010-2013) up to become rows 2-5.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2:A7").EntireRow.Delete()
